$d = $word.ActiveDocument

# Locate the "K1." paragraph in section "K. Overig" by its distinctive
# trailing question text, so the new paragraph can be inserted right
# after it (and right before the "K2." paragraph).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*is het saldo juridisch vastgesteld?*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the K1 paragraph (saldo-overzichten question)."
}

# Create a brand-new empty paragraph directly after the K1 paragraph.
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()
$insertRange = $newPara.Range

# Fill the new (still-empty) paragraph with the "Feitelijk: " bold label
# run followed by the italic factual-answer run, using the exact OOXML
# so formatting/run-splitting matches precisely.
$answer = "Nee, Freca heeft geen jaarlijkse saldo-overzichten verstrekt. " +
          "Dinck heeft geen bezwaar gemaakt. Er is geen juridisch " +
          "vastgesteld saldo " + [char]0x2014 + " het huidige saldo komt " +
          "uitsluitend uit de Snelstart-boekhouding."

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
  '<w:pPr>' +
    '<w:spacing w:after="80"/>' +
    '<w:ind w:left="360"/>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>' +
      '<w:b/>' +
      '<w:bCs/>' +
      '<w:sz w:val="22"/>' +
      '<w:szCs w:val="22"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve">Feitelijk: </w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>' +
      '<w:i/>' +
      '<w:iCs/>' +
      '<w:sz w:val="22"/>' +
      '<w:szCs w:val="22"/>' +
    '</w:rPr>' +
    "<w:t xml:space=`"preserve`">$answer</w:t>" +
  '</w:r>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xmlFragment)
